# Auto-generated edit script: updates crypto price/volume table cells
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number-format on the specific Price cells whose new values
# would otherwise be auto-parsed as numeric by Excel (losing the original
# "text" cell type used throughout column D).
$textCells = @("D5","D6","D8","D9","D11","D15","D16","D18","D19","D20","D23","D24","D26","D29","D32","D33","D35","D36","D37","D38","D40","D42","D43","D45","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values from the commit.
$ws.Range("D2").Value = "58.333.78"
$ws.Range("E2").Value = "  -4.28%  "
$ws.Range("D3").Value = "2.611.48"
$ws.Range("E3").Value = "  -4.13%  "
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").Value = "520.00"
$ws.Range("E5").Value = "  -1.77%  "
$ws.Range("D6").Value = "141.84"
$ws.Range("E6").Value = "  -3.69%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").Value = "0.567"
$ws.Range("E8").Value = "  -2.38%  "
$ws.Range("D9").Value = "6.54"
$ws.Range("E9").Value = "  -8.69%  "
$ws.Range("E10").Value = "  -3.72%  "
$ws.Range("D11").Value = "0.335"
$ws.Range("E11").Value = "  -2.06%  "
$ws.Range("E12").Value = "  +0.82%  "
$ws.Range("D13").Value = "3.067.95"
$ws.Range("E13").Value = "  -4.60%  "
$ws.Range("D14").Value = "58.286.43"
$ws.Range("E14").Value = "  -4.41%  "
$ws.Range("D15").Value = "20.83"
$ws.Range("E15").Value = "  -3.49%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "0.0000135"
$ws.Range("E16").Value = "  -2.65%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.614.24"
$ws.Range("E17").Value = "  -6.43%  "
$ws.Range("D18").Value = "336.51"
$ws.Range("E18").Value = "  -2.98%  "
$ws.Range("D19").Value = "4.39"
$ws.Range("E19").Value = "  -3.11%  "
$ws.Range("D20").Value = "10.34"
$ws.Range("E20").Value = "  -2.68%  "
$ws.Range("E21").Value = "  -3.36%  "
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("D23").Value = "64.85"
$ws.Range("E23").Value = "  +1.81%  "
$ws.Range("D24").Value = "0.413"
$ws.Range("E24").Value = "  -2.02%  "
$ws.Range("E25").Value = "  -3.61%  "
$ws.Range("D26").Value = "0.997"
$ws.Range("E26").Value = "  +0.16%  "
$ws.Range("E27").Value = "  -3.61%  "
$ws.Range("D28").Value = "0.0₃0784"
$ws.Range("E28").Value = "  -5.46%  "
$ws.Range("D29").Value = "6.52"
$ws.Range("E29").Value = "  -3.89%  "
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("E31").Value = "  -1.46%  "
$ws.Range("D32").Value = "18.71"
$ws.Range("E32").Value = "  -2.55%  "
$ws.Range("D33").Value = "150.10"
$ws.Range("E33").Value = "  -0.10%  "
$ws.Range("E34").Value = "  -4.99%  "
$ws.Range("D35").Value = "1.18"
$ws.Range("E35").Value = "  -4.77%  "
$ws.Range("D36").Value = "0.887"
$ws.Range("E36").Value = "  -4.20%  "
$ws.Range("D37").Value = "0.846"
$ws.Range("E37").Value = "  -6.73%  "
$ws.Range("D38").Value = "36.15"
$ws.Range("E38").Value = "  -2.97%  "
$ws.Range("E39").Value = "  -7.57%  "
$ws.Range("D40").Value = "3.61"
$ws.Range("E40").Value = "  -2.46%  "
$ws.Range("E41").Value = "  +0.28%  "
$ws.Range("D42").Value = "0.600"
$ws.Range("E42").Value = "  -4.07%  "
$ws.Range("D43").Value = "0.0968"
$ws.Range("E43").Value = "  -2.43%  "
$ws.Range("E44").Value = "  +0.98%  "
$ws.Range("D45").Value = "266.51"
$ws.Range("E45").Value = "  -6.17%  "
$ws.Range("E46").Value = "  -7.07%  "
$ws.Range("E47").Value = "  -3.38%  "
$ws.Range("D48").Value = "2.020.51"
$ws.Range("E48").Value = "  -5.05%  "
$ws.Range("E49").Value = "  -2.46%  "
$ws.Range("E50").Value = "  -8.55%  "
$ws.Range("D51").Value = "18.14"
$ws.Range("E51").Value = "  -6.96%  "
